$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the TEST column (E) to fit the new test-case text
$ws.Columns.Item(5).ColumnWidth = 26.02

# Fill in the new "TEST" column values (Test Cases) for each requirement row
$ws.Range("E4").Value = "TC_3, "
$ws.Range("E5").Value = "TC_7, TC_8, TC_9"
$ws.Range("E6").Value = "TC_1, TC_2`nTC_4, TC_5, TC_6"
$ws.Range("E7").Value = "TC_10"

# Row 6's test entry wraps onto two lines -- enable wrap and grow the row
$ws.Range("E6").WrapText = $true
$ws.Range("E6").EntireRow.RowHeight = 30

# CODE column was blank for the "NA" row -- mirror the CRS/SRS/DESIGN "NA"
$ws.Range("D8").Value = "NA"

# Update the view: scroll right slightly and reselect where the user left off
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E15").Select()
